$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of MFE fly data to append starting at row 234 (vial, treatment, time_hours, males, females)
$newData = @(
  @(1, "conditioned", 380, 2, 4),
  @(1, "unconditioned", 380, 0, 0),
  @(2, "conditioned", 380, 0, 0),
  @(2, "unconditioned", 380, 0, 0),
  @(3, "conditioned", 380, 0, 0),
  @(3, "unconditioned", 380, 0, 0),
  @(4, "conditioned", 380, 0, 0),
  @(4, "unconditioned", 380, 2, 1),
  @(5, "conditioned", 380, 0, 0),
  @(5, "unconditioned", 380, 0, 1),
  @(6, "conditioned", 380, 2, 2),
  @(6, "unconditioned", 380, 1, 1),
  @(7, "conditioned", 380, 0, 0),
  @(7, "unconditioned", 380, 3, 0),
  @(8, "conditioned", 380, 0, 1),
  @(8, "unconditioned", 380, 1, 1),
  @(9, "conditioned", 380, 0, 2),
  @(9, "unconditioned", 380, 0, 1),
  @(10, "conditioned", 380, 0, 0),
  @(10, "unconditioned", 380, 1, 1),
  @(11, "conditioned", 380, 0, 3),
  @(11, "unconditioned", 380, 0, 0),
  @(12, "conditioned", 380, 0, 1),
  @(12, "unconditioned", 380, 0, 0),
  @(13, "conditioned", 380, 0, 0),
  @(13, "unconditioned", 380, 1, 1),
  @(14, "conditioned", 380, 2, 3),
  @(14, "unconditioned", 380, 0, 0),
  @(15, "unconditioned", 380, 0, 0),
  @(1, "conditioned", 389, 1, 2),
  @(1, "unconditioned", 389, 0, 0),
  @(2, "conditioned", 389, 0, 1),
  @(2, "unconditioned", 389, 0, 0),
  @(3, "conditioned", 389, 0, 0),
  @(3, "unconditioned", 389, 0, 0),
  @(4, "conditioned", 389, 2, 2),
  @(4, "unconditioned", 389, 1, 5),
  @(5, "conditioned", 389, 0, 0),
  @(5, "unconditioned", 389, 2, 1),
  @(6, "conditioned", 389, 1, 3),
  @(6, "unconditioned", 389, 1, 0),
  @(7, "conditioned", 389, 0, 0),
  @(7, "unconditioned", 389, 1, 1),
  @(8, "conditioned", 389, 2, 2),
  @(8, "unconditioned", 389, 0, 0),
  @(9, "conditioned", 389, 1, 5),
  @(9, "unconditioned", 389, 0, 0),
  @(10, "conditioned", 389, 1, 1),
  @(10, "unconditioned", 389, 0, 2),
  @(11, "conditioned", 389, 0, 2),
  @(11, "unconditioned", 389, 0, 0),
  @(12, "conditioned", 389, 3, 4),
  @(12, "unconditioned", 389, 0, 0),
  @(13, "conditioned", 389, 0, 0),
  @(13, "unconditioned", 389, 1, 1),
  @(14, "conditioned", 389, 2, 1),
  @(14, "unconditioned", 389, 0, 0),
  @(15, "unconditioned", 389, 0, 0),
  @(1, "conditioned", 404, 1, 2),
  @(1, "unconditioned", 404, 0, 0),
  @(2, "conditioned", 404, 0, 0),
  @(2, "unconditioned", 404, 0, 0),
  @(3, "conditioned", 404, 0, 0),
  @(3, "unconditioned", 404, 0, 0),
  @(4, "conditioned", 404, 0, 0),
  @(4, "unconditioned", 404, 0, 1),
  @(5, "conditioned", 404, 0, 0),
  @(5, "unconditioned", 404, 0, 0),
  @(6, "conditioned", 404, 0, 0),
  @(6, "unconditioned", 404, 0, 0),
  @(7, "conditioned", 404, 0, 0),
  @(7, "unconditioned", 404, 0, 0),
  @(8, "conditioned", 404, 1, 2),
  @(8, "unconditioned", 404, 0, 1),
  @(9, "conditioned", 404, 0, 2),
  @(9, "unconditioned", 404, 0, 0),
  @(10, "conditioned", 404, 0, 3),
  @(10, "unconditioned", 404, 1, 1),
  @(11, "conditioned", 404, 1, 0),
  @(11, "unconditioned", 404, 0, 0),
  @(12, "conditioned", 404, 2, 1),
  @(12, "unconditioned", 404, 0, 0),
  @(13, "conditioned", 404, 0, 0),
  @(13, "unconditioned", 404, 0, 0),
  @(14, "conditioned", 404, 0, 0),
  @(14, "unconditioned", 404, 0, 0),
  @(15, "unconditioned", 404, 0, 0)
)

$startRow = 234
$rowCount = $newData.Count
$colCount = 5

$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i, $j] = $newData[$i][$j]
    }
}

$endRow = $startRow + $rowCount - 1
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 5))
$targetRange.Value = $arr

$ws.Range("D293").Select()
